# formula: support for more statistical functions
#   AVERAGE, AVERAGEA, COUNT, COUNTA, COUNTBLANK, MIN, MAX, MEDIAN
#
# Adds a new "Statistics" worksheet (after "Indexing") exercising the new
# statistical functions, matching the shape of the other function-coverage
# sheets already in this workbook (Logical / Math and Trig / Text / Indexing).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheet at the end of the tab strip and make it the active one
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Statistics"
$ws.Activate()

# ---------------------------------------------------------------------------
# 2. Header row (bold) - one column per new function
# ---------------------------------------------------------------------------
$headers = @(
    @("E1", "AVERAGE"),
    @("F1", "AVERAGEA"),
    @("G1", "COUNT"),
    @("H1", "COUNTA"),
    @("I1", "COUNTBLANK"),
    @("J1", "MIN"),
    @("K1", "MAX"),
    @("L1", "MEDIAN")
)
foreach ($h in $headers) {
    $ws.Range($h[0]).Value = $h[1]
    $ws.Range($h[0]).Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 3. Raw sample data (A2:C4 numeric block, B5/C5 extra row, A6/B6 booleans)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 6
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 7
$ws.Range("B5").Value = "asd"
$ws.Range("C5").Value = 12

$ws.Range("A6").Formula = "=TRUE()"
$ws.Range("A6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("B6").Formula = "=FALSE()"
$ws.Range("B6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# ---------------------------------------------------------------------------
# 4. Formulas exercising each new statistical function
# ---------------------------------------------------------------------------
$formulas = @(
    @("E2", "AVERAGE()"),
    @("F2", "AVERAGEA(A2:C4)"),
    @("G2", "COUNT()"),
    @("H2", "COUNTA()"),
    @("I2", "COUNTBLANK()"),
    @("J2", "MIN()"),
    @("K2", "MAX()"),
    @("L2", "MEDIAN(A2:C2)"),

    @("E3", "AVERAGE(A2:C2)"),
    @("F3", "AVERAGEA(0)"),
    @("G3", "COUNT(A2:C2)"),
    @("H3", "COUNTA(A2:C2)"),
    @("I3", "COUNTBLANK(A3:C4)"),
    @("J3", "MIN(A2:C4)"),
    @("K3", "MAX(A3:C4)"),
    @("L3", "MEDIAN()"),

    @("E4", "AVERAGE(C2)"),
    @("F4", "AVERAGEA()"),
    @("G4", "COUNT(A3:C4)"),
    @("H4", "COUNTA(A4:C5)"),
    @("I4", "COUNTBLANK(A4:C6)"),
    @("J4", "MIN(B3:C4)"),
    @("K4", "MAX(A2:B3)"),
    @("L4", "MEDIAN(A4)"),

    @("E5", "AVERAGE(A2:A4)"),
    @("F5", "AVERAGEA(A5:C5)"),
    @("G5", "COUNT(A5:C5)"),
    @("H5", "COUNTA(A5:C6)"),
    @("I5", "COUNTBLANK(A7:D9)"),
    @("J5", "MIN(J6:J9)"),
    @("K5", "MAX(K6:K9)"),
    @("L5", "MEDIAN(A3:C4)"),

    @("E6", "AVERAGE(A2:A5)"),
    @("F6", "AVERAGE(B5:C6)"),
    @("G6", "COUNT(A6:C6)"),
    @("H6", "COUNTA(A7:D9)"),
    @("L6", "MEDIAN(A2:C6)"),

    @("E7", "AVERAGE(C2:C5)"),
    @("F7", "AVERAGEA(A5:C6)"),
    @("G7", "COUNT(A1:C6)"),

    @("E8", "AVERAGE(B2:B5)"),
    @("E9", "AVERAGE(B2:B4)"),
    @("E10", "AVERAGE(A2:A6)"),
    @("E11", "AVERAGE(A6:B6)"),
    @("E12", "AVERAGE(A6:C6)")
)
foreach ($f in $formulas) {
    $ws.Range($f[0]).Formula = "=" + $f[1]
}

# ---------------------------------------------------------------------------
# 5. Column widths on the new sheet (best-effort match of authoring tool's
#    autofit widths)
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").EntireColumn.ColumnWidth = 7.67
$ws.Range("E1").EntireColumn.ColumnWidth = 11.86
$ws.Range("F1").EntireColumn.ColumnWidth = 14.42
$ws.Range("G1").EntireColumn.ColumnWidth = 7.67
$ws.Range("H1").EntireColumn.ColumnWidth = 8.35
$ws.Range("I1").EntireColumn.ColumnWidth = 12.4

# ---------------------------------------------------------------------------
# 6. Minor column-width touch-ups on pre-existing sheets (side effect of the
#    original authoring tool re-flowing autofit widths once the new sheet's
#    strings were added).
# ---------------------------------------------------------------------------
$wsText = $wb.Worksheets.Item("Text")
$wsText.Range("E1").EntireColumn.ColumnWidth = 9.02
$wsText.Range("F1").EntireColumn.ColumnWidth = 7.53
$wsText.Range("O1").EntireColumn.ColumnWidth = 8.75
$wsText.Range("P1").EntireColumn.ColumnWidth = 8.48

$wsIndexing = $wb.Worksheets.Item("Indexing")
$wsIndexing.Range("G1").EntireColumn.ColumnWidth = 9.02
$wsIndexing.Range("H1").EntireColumn.ColumnWidth = 8.62

# Re-activate the new sheet/cell selection to match the authoring tool's
# final cursor position.
$ws.Activate()
$ws.Range("L7").Select()
